$wb = $excel.ActiveWorkbook

# The workbook originally opened with "Absolute History" (sheet5) as the
# active/selected tab. Walk through the history sheets adding a new
# "Comments" header in column E, finishing on "Withdraw History" so that
# it ends up as the active tab (matching the saved workbook view).

# --- Absolute History (was the active sheet) ---
$wsAbsolute = $wb.Worksheets.Item("Absolute History")
$wsAbsolute.Activate()
$wsAbsolute.Range("E1").Value = "Comments"
[void]$wsAbsolute.Range("E5").Select()

# --- Transfer History ---
$wsTransfer = $wb.Worksheets.Item("Transfer History")
$wsTransfer.Activate()
$wsTransfer.Range("E1").Value = "Comments"
[void]$wsTransfer.Range("E1").Select()

# --- Deposit History ---
$wsDeposit = $wb.Worksheets.Item("Deposit History")
$wsDeposit.Activate()
$wsDeposit.Range("E1").Value = "Comments"
[void]$wsDeposit.Range("E1").Select()

# --- Withdraw History (ends up the active tab) ---
$wsWithdraw = $wb.Worksheets.Item("Withdraw History")
$wsWithdraw.Activate()
$wsWithdraw.Range("E1").Value = "Comments"
[void]$wsWithdraw.Range("E1").Select()
